$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (row 4) with David Bejenariu's data
$ws.Range("A4").Value = "David Bejenariu"
$ws.Range("B4").Value = "david.bejenariu@gmail.com"
$ws.Range("C4").Value = "Black"
$ws.Range("D4").Value = "Relevance"

# Turn the email into a mailto hyperlink, matching the style used by the
# other email cells (B2, B3)
$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:david.bejenariu@gmail.com")
$ws.Range("B4").Style = "Hyperlink"

# Update the active selection like the author's session ended up
$ws.Range("C6").Select() | Out-Null
